$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I26").Value = "development front-end af"
$ws.Range("J26").Value = "1,5 uur"
$ws.Range("K26").Value = (Get-Date -Year 2017 -Month 10 -Day 23)
$ws.Range("L26").Value = "Het ging redelijk simpel op een paar bugs na"
Write-Host "done"
